$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the old "celkem" summary block (old rows 26-29),
# shifting that block down to rows 30-33 to make room for new field-visit data.
$ws.Rows("26:29").Insert()

# New data rows recorded on the "last field visit" (14.8.2023), matching the
# existing "mapování v terénu" / "posprocessing" entries used elsewhere.
$ws.Range("C24").Value = 45151
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = "mapování v terénu"
$ws.Range("F24").Value = 112

$ws.Range("C25").Value = 45152
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = "mapování v terénu"
$ws.Range("F25").Value = 112

$ws.Range("C26").Value = 45152
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = "posprocessing"

# Give the new date cells (C24:C26) the same date display format as the rest
# of column C, and carry that same formatting one row further (C27), with no
# value, matching the leftover formatted-but-empty cell below the data.
$ws.Range("C23").Copy()
$ws.Range("C24:C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-point the totals formulas (now at rows 31-33 after the insert above) at
# the expanded data range.
$ws.Range("D31").Formula = "=SUM(D6:D29)"
$ws.Range("D32").Formula = "=D30*D31"
$ws.Range("D33").Formula = "=D32+SUM(F6:F29)+SUM(G6:G29)"

# Restore the selection left by the last edit.
$ws.Range("D32").Select()
